# Apply updated cryptocurrency price/volume figures to Sheet1.
# Values are stored as text (matching the workbook's existing inlineStr
# convention), so each target cell is explicitly formatted as Text before
# the new value is written -- this prevents Excel from re-interpreting
# strings such as "6.740" or "0.9100" as numbers and silently dropping
# the significant trailing zero.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    'D2' = '265.62'; 'E2' = '1.78%'
    'D3' = '26.62'; 'E3' = '-1.51%'
    'D4' = '4.698'; 'E4' = '-0.12%'
    'E5' = '-1.56%'
    'D6' = '6.740'; 'E6' = '0.80%'
    'D7' = '0.8515'; 'E7' = '0.17%'
    'D8' = '0.9100'; 'E8' = '-0.75%'
    'E9' = '0.15%'
    'D10' = '0.05151'; 'E10' = '10.25%'
    'D11' = '0.07102'; 'E11' = '0.20%'
    'D12' = '0.03147'; 'E12' = '-0.33%'
    'D13' = '0.09031'; 'E13' = '-0.09%'
    'D14' = '0.001528'; 'E14' = '0.20%'
    'D15' = '0.0006056'; 'E15' = '-1.94%'
    'D16' = '0.005936'; 'E16' = '-2.96%'
    'D17' = '3.451'; 'E17' = '0.00%'
    'E18' = '-0.02%'
    'E19' = '0.29%'
    'E21' = '-1.41%'
    'D22' = '4.123'; 'E22' = '1.04%'
    'D23' = '0.04242'; 'E23' = '0.21%'
    'D24' = '0.001179'; 'E24' = '-3.15%'
    'D25' = '0.004062'; 'E25' = '6.85%'
    'E26' = '0.02%'
    'E27' = '6.58%'
    'D40' = '0.03924'; 'E40' = '1.00%'
    'D41' = '0.1116'; 'E41' = '0.49%'
    'D42' = '0.004196'; 'E42' = '2.47%'
    'E43' = '-3.57%'
    'D44' = '0.01149'; 'E44' = '-29.59%'
    'D45' = '0.00005089'; 'E45' = '-1.41%'
    'D46' = '0.00000000751'; 'E46' = '0.03%'
    'D48' = '0.2582'; 'E48' = '54.95%'
    'E49' = '0.03%'
    'E50' = '0.03%'
}

foreach ($ref in $updates.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$ref]
}
